# Insert a new record row at row 199 (pushing the existing rows 199-243 down to 200-244)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 199; existing row 199 (and below) shift down to 200.
$ws.Rows.Item(199).Insert()

# Fill in the new row 199 with the observation values.
$ws.Cells.Item(199, 1).Value = 3
$ws.Cells.Item(199, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(199, 3).Value = "Coquimbo"
$ws.Cells.Item(199, 4).Value = 44511
$ws.Cells.Item(199, 5).Value = 5
$ws.Cells.Item(199, 6).Value = 100112031
$ws.Cells.Item(199, 7).Value = "Poroto verde"
$ws.Cells.Item(199, 8).Value = "Magnum"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 73
$ws.Cells.Item(199, 11).Value = 42000
$ws.Cells.Item(199, 12).Value = 43000
$ws.Cells.Item(199, 13).Value = 42521
$ws.Cells.Item(199, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(199, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(199, 16).Value = 1701
$ws.Cells.Item(199, 17).Value = 25
$ws.Cells.Item(199, 18).Value = "Hortaliza"
